$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 2 (done first): "if the child has any special need|s or
# circumstances." text.
#   Remove the existing "_GoBack" bookmark that splits "need" from
#   "s or circumstances." and merge those two runs into one run (leaving the
#   preceding "At the hearing, you can tell the court " run untouched).
#
# This has to happen before Change 1 below because Word only allows a single
# bookmark per name - re-adding "_GoBack" at the <<courtName>> location
# would otherwise just relocate this existing one instead of creating an
# independent bookmark.
# ---------------------------------------------------------------------------

$needBm = $d.Bookmarks("_GoBack")
$needBmEnd = $needBm.End
$needParaEnd = $needBm.Range.Paragraphs(1).Range.End
$needBm.Delete()

# Re-assigning the text of the run that follows the (now deleted) bookmark,
# using a range that reaches all the way to the paragraph's end, triggers
# Word's adjacent-run coalescing for same-formatted runs once the bookmark
# "wall" between them is gone - merging "s or circumstances. " together with
# the preceding "if the child has any special need" run, without touching
# the earlier, differently-edited "At the hearing..." run.
$afterBookmark = $d.Range($needBmEnd, $needParaEnd)
$afterBookmark.Text = "s or circumstances. "

# ---------------------------------------------------------------------------
# Change 1: "<<courtName>>" paragraph.
#   Split the single run "courtName>>" into three runs: "court", "Name", ">>"
#   and insert the "_GoBack" bookmark between "Name" and ">>".
# ---------------------------------------------------------------------------

# Locate the "<<courtName>>" paragraph and compute absolute character offsets
# for the boundaries we need to split on ("court"|"Name" and "Name"|">>").
$courtPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^<<courtName>>") {
        $courtPara = $p
        break
    }
}

$fullText = $courtPara.Range.Text
$paraStart = $courtPara.Range.Start
$courtNameOffset = $fullText.IndexOf("courtName")
$courtSplit = $paraStart + $courtNameOffset + "court".Length
$nameSplit = $paraStart + $courtNameOffset + "courtName".Length

# Split "court" from "Name" by toggling (and restoring) Bold across "court" -
# this forces the run boundary without altering the visible formatting.
$courtRun = $d.Range($paraStart + $courtNameOffset, $courtSplit)
$courtRun.Font.Bold = 1
$courtRun.Font.Bold = 0

# Insert the "_GoBack" bookmark right after "Name" (before ">>"); adding a
# collapsed bookmark at this point also splits "Name" from ">>" into
# separate runs.
$bookmarkPoint = $d.Range($nameSplit, $nameSplit)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
